$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=9;   C=69579;  E=191683877},
    @{Row=10;  C=278214; E=1752681420},
    @{Row=17;  C=134754; E=296820206},
    @{Row=38;  C=27073;  E=83546230},
    @{Row=54;  C=17837;  E=32098865},
    @{Row=69;  C=20737;  E=62181544},
    @{Row=99;  C=136581; E=863240545},
    @{Row=126; C=5647;   E=8181370},
    @{Row=169; C=562674; E=1286269320},
    @{Row=170; C=367599; E=2848423734},
    @{Row=174; C=357393; E=1020291034},
    @{Row=175; C=125703; E=815960902},
    @{Row=179; C=235816; E=813759583},
    @{Row=203; C=13108;  E=33026608},
    @{Row=205; C=11135;  E=44608881},
    @{Row=243; C=28204;  E=90602697},
    @{Row=257; C=182554; E=1063849900},
    @{Row=262; C=38988;  E=124760356},
    @{Row=266; C=71667;  E=219458739},
    @{Row=311; C=190864; E=586829007},
    @{Row=323; C=94726;  E=178876673}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
